$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 819.2727
$ws.Range("I19").Value = 250
$ws.Range("J19").Value = 945.7778
$ws.Range("K19").Value = 250
$ws.Range("L19").Value = 945.7778
$ws.Range("M19").Value = -75
$ws.Range("N19").Value = -1295.7778
$ws.Range("H41").Value = 375.75
$ws.Range("I41").Value = 292.66666
$ws.Range("K41").Value = 292.66666
$ws.Range("M41").Value = 147.33334
$ws.Range("H53").Value = 624
$ws.Range("I53").Value = 624
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 624
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 13
$ws.Range("N53").Value = $null
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = $null
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = $null
$ws.Range("H106").Value = 1500
$ws.Range("I106").Value = 1500
$ws.Range("K106").Value = 1500
$ws.Range("M106").Value = -869
$ws.Range("H137").Value = 3467.3
$ws.Range("I137").Value = 3162.5
$ws.Range("K137").Value = 9487.5
$ws.Range("M137").Value = -6937.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1555.1111
$ws.Range("I61").Value = 1471.5714
$ws.Range("K61").Value = 1471.5714
$ws.Range("M61").Value = -1259.5714
$ws.Range("H74").Value = 3644.9
$ws.Range("I74").Value = 3644.9
$ws.Range("K74").Value = 3644.9
$ws.Range("M74").Value = -2770.9
$ws.Range("H77").Value = 3644.9
$ws.Range("I77").Value = 3644.9
$ws.Range("K77").Value = 18224.5
$ws.Range("M77").Value = -13856.5
$ws.Range("H107").Value = 75000
$ws.Range("J107").Value = 75000
$ws.Range("L107").Value = 75000
$ws.Range("N107").Value = -82680
$ws.Range("H136").Value = 1555.1111
$ws.Range("I136").Value = 1471.5714
$ws.Range("K136").Value = 4414.7142
$ws.Range("M136").Value = -1864.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 9333
$ws.Range("J19").Value = 9999.5
$ws.Range("L19").Value = 9999.5
$ws.Range("N19").Value = -10345.5
$ws.Range("H80").Value = 1520
$ws.Range("I80").Value = 1400
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 1400
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -402
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 1520
$ws.Range("I83").Value = 1400
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 7000
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -2008
$ws.Range("N83").Value = -19984
$ws.Range("H134").Value = 2442.7896
$ws.Range("I134").Value = 2442.7896
$ws.Range("K134").Value = 7328.3688
$ws.Range("M134").Value = -4793.3688

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 815
$ws.Range("I22").Value = 125
$ws.Range("J22").Value = 1850
$ws.Range("K22").Value = 125
$ws.Range("L22").Value = 1850
$ws.Range("M22").Value = 225
$ws.Range("N22").Value = -2550
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null
$ws.Range("H134").Value = 500
$ws.Range("I134").Value = 500
$ws.Range("K134").Value = 1500
$ws.Range("M134").Value = 1035
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 135.25
$ws.Range("I60").Value = 97
$ws.Range("K60").Value = 291
$ws.Range("M60").Value = -40
$ws.Range("H114").Value = 1382.75
$ws.Range("I114").Value = 1500
$ws.Range("J114").Value = 1031
$ws.Range("K114").Value = 4500
$ws.Range("L114").Value = 3093
$ws.Range("M114").Value = -1246
$ws.Range("N114").Value = -9601
$ws.Range("H137").Value = 6165.6
$ws.Range("J137").Value = 4809.3335
$ws.Range("L137").Value = 14428.0005
$ws.Range("N137").Value = -24628.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 24000
$ws.Range("J101").Value = 24000
$ws.Range("L101").Value = 24000
$ws.Range("N101").Value = -30490
$ws.Range("H132").Value = 1165.3334
$ws.Range("I132").Value = 1165.3334
$ws.Range("K132").Value = 3496.0002
$ws.Range("M132").Value = -966.0001999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 2153.6
$ws.Range("I14").Value = 966.3333
$ws.Range("K14").Value = 966.3333
$ws.Range("M14").Value = -794.3333
$ws.Range("H38").Value = 95033
$ws.Range("J38").Value = 95033
$ws.Range("L38").Value = 95033
$ws.Range("N38").Value = -95853
$ws.Range("H101").Value = 31262
$ws.Range("J101").Value = 31262
$ws.Range("L101").Value = 31262
$ws.Range("N101").Value = -37752
$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -8100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 158851
$ws.Range("J103").Value = 158851
$ws.Range("L103").Value = 158851
$ws.Range("N103").Value = -161195
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988
$ws.Range("H117").Value = 18350.5
$ws.Range("J117").Value = 18350.5
$ws.Range("L117").Value = 18350.5
$ws.Range("N117").Value = -27528.5
$ws.Range("H126").Value = 1413.8334
$ws.Range("I126").Value = 1387.8182
$ws.Range("K126").Value = 4163.4546
$ws.Range("M126").Value = -1693.4546
$ws.Range("H132").Value = 3242.7144
$ws.Range("I132").Value = 2674.75
$ws.Range("K132").Value = 8024.25
$ws.Range("M132").Value = -5494.25
$ws.Range("H136").Value = 7922.222
$ws.Range("I136").Value = 6946.231
$ws.Range("J136").Value = 10459.8
$ws.Range("K136").Value = 20838.693
$ws.Range("L136").Value = 31379.4
$ws.Range("M136").Value = -18288.693
$ws.Range("N136").Value = -36479.39999999999

Write-Host "Applied all Golem_Profits updates"